$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Resize the 4 columns of the (only) table in the document so the
#    underlying <w:tblGrid> gridCol widths become 1628, 2083, 2959,
#    2539 (dxa) i.e. 81.4, 104.15, 147.95, 126.95 points.
# -----------------------------------------------------------------
$t = $d.Tables.Item(1)
$t.Columns.Item(1).Width = 81.4
$t.Columns.Item(2).Width = 104.15
$t.Columns.Item(3).Width = 147.95
$t.Columns.Item(4).Width = 126.95

# -----------------------------------------------------------------
# 2) Fix the "endeudamient_indicador_dp" typo -> "endeudamiento_indicador_dp"
#    (missing "o") in the "Razon de endeudamiento" row, 3rd column.
#    The merge-then-resplit sequence below reproduces the exact run
#    layout Word itself produces when a user clicks right after
#    "endeudamient" and types the missing "o": the previously separate
#    "{", "endeudamient_indicador_dp" and "}" runs (with stale
#    spell-check proofErr markers around the middle one) collapse into
#    a single run carrying the corrected text, which is then re-split
#    into "{endeudamient" / "o" / "_indicador_dp}" at the edit point.
# -----------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("{endeudamient_indicador_dp}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p0 = $rng.Start
$rng.Text = "{endeudamiento_indicador_dp}"

# Re-split the merged run into three runs at the point where the "o"
# was inserted (offsets 13 and 14 from the start of "{endeudamiento_indicador_dp}").
$mid = $d.Range($p0 + 13, $p0 + 14)
$mid.Font.Bold = 1
$mid.Font.Bold = 0
